# Generate Report for Handoff
# Adds two newly-discovered localization files
#   6dbee821-dc55-4aa6-b9f8-a85a8a7c608a  (sorted before 8d18e275...)
#   d1c0fc93-1653-4056-9649-c442a3e88578  (sorted after  8d18e275...)
# to the Overview/zh-cn/de-de report sheets, each "Ready for handoff".

$wb = $excel.ActiveWorkbook

function Add-ReportHyperlink($ws, $cellRef, $target, $display) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Make room: row 3 becomes the new "6dbee821" row (old row3 -> row4), then
# row 5 becomes the new "d1c0fc93" row (old row4, now holding 8d18e275,
# stays put at row4).
$wsOv.Range("A3").EntireRow.Insert()
$wsOv.Range("A5").EntireRow.Insert()

$wsOv.Range("A3").Value = "6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.md"
$wsOv.Range("B3").Value = "Ready for handoff"
$wsOv.Range("C3").Value = "Ready for handoff"
$wsOv.Range("D3").Value = "2016-31-19 08:31:17"

$wsOv.Range("A5").Value = "d1c0fc93-1653-4056-9649-c442a3e88578.md"
$wsOv.Range("B5").Value = "Ready for handoff"
$wsOv.Range("C5").Value = "Ready for handoff"
$wsOv.Range("D5").Value = "2016-31-19 08:31:17"

# Hyperlinks collection doesn't follow row-inserts, so rebuild all of them
# (ranges + display text) for the "File Name" column in final order.
$wsOv.Hyperlinks.Delete()
Add-ReportHyperlink $wsOv "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f3f0bcc76a05ce10f5c256e5d9009383cc543374/e2e/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md" "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md"
Add-ReportHyperlink $wsOv "A3" "https://github.com/OpenLocalizationTest/oltest/blob/ed81bff2ae9620a8cb3ed916d5d61fbfca77ec35/e2e/6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.md" "6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.md"
Add-ReportHyperlink $wsOv "A4" "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/8d18e275-99a5-4481-83b7-9dcffb239eac.md" "8d18e275-99a5-4481-83b7-9dcffb239eac.md"
Add-ReportHyperlink $wsOv "A5" "https://github.com/OpenLocalizationTest/oltest/blob/b87ec7a0b392ebf3a46797c41e6174f22f049a84/e2e/d1c0fc93-1653-4056-9649-c442a3e88578.md" "d1c0fc93-1653-4056-9649-c442a3e88578.md"

# ---------------------------------------------------------------------------
# Language detail sheets: "zh-cn" and "de-de" share the same layout:
# Source File Name | File Extension | Status | Latest Handoff File |
# Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Handoff Reason | Dependency From | Error Detail
# ---------------------------------------------------------------------------

$langs = @(
  @{
    Sheet = "zh-cn"
    Suffix = "zh-cn"
    Row3 = @{ D = "6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.ed81bff2ae9620a8cb3ed916d5d61fbfca77ec35.zh-cn.xlf"; E = "2016-03-19 08:31:14" }
    Row4 = @{ D = "8d18e275-99a5-4481-83b7-9dcffb239eac.38bef17d5ca5f7108f12f6ff593a277e5c1978d7.zh-cn.xlf"; E = "2016-03-19 08:30:02" }
    Row5 = @{ D = "d1c0fc93-1653-4056-9649-c442a3e88578.b87ec7a0b392ebf3a46797c41e6174f22f049a84.zh-cn.xlf"; E = "2016-03-19 08:31:14" }
  },
  @{
    Sheet = "de-de"
    Suffix = "de-de"
    Row3 = @{ D = "6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.ed81bff2ae9620a8cb3ed916d5d61fbfca77ec35.de-de.xlf"; E = "2016-03-19 08:31:17" }
    Row4 = @{ D = "8d18e275-99a5-4481-83b7-9dcffb239eac.38bef17d5ca5f7108f12f6ff593a277e5c1978d7.de-de.xlf"; E = "2016-03-19 08:30:08" }
    Row5 = @{ D = "d1c0fc93-1653-4056-9649-c442a3e88578.b87ec7a0b392ebf3a46797c41e6174f22f049a84.de-de.xlf"; E = "2016-03-19 08:31:17" }
  }
)

foreach ($lang in $langs) {
    $ws = $wb.Worksheets.Item($lang.Sheet)
    $suffix = $lang.Suffix

    # Same row-shuffle as the Overview sheet.
    $ws.Range("A3").EntireRow.Insert()
    $ws.Range("A5").EntireRow.Insert()

    # --- new row 3: 6dbee821... ---
    $ws.Range("A3").Value = "6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.md"
    $ws.Range("B3").Value = ".md"
    $ws.Range("C3").Value = "Ready for handoff"
    $ws.Range("D3").Value = $lang.Row3.D
    $ws.Range("E3").Value = $lang.Row3.E
    $ws.Range("H3").Value = "0001-01-01 00:00:00"
    $ws.Range("I3").Value = "Include"

    # --- row 4: 8d18e275... (pre-existing row, now shifted down) ---
    $ws.Range("A4").Value = "8d18e275-99a5-4481-83b7-9dcffb239eac.md"
    $ws.Range("B4").Value = ".md"
    $ws.Range("C4").Value = "Ready for handoff"
    $ws.Range("D4").Value = $lang.Row4.D
    $ws.Range("E4").Value = $lang.Row4.E
    $ws.Range("H4").Value = "0001-01-01 00:00:00"
    $ws.Range("I4").Value = "Include"

    # --- new row 5: d1c0fc93... ---
    $ws.Range("A5").Value = "d1c0fc93-1653-4056-9649-c442a3e88578.md"
    $ws.Range("B5").Value = ".md"
    $ws.Range("C5").Value = "Ready for handoff"
    $ws.Range("D5").Value = $lang.Row5.D
    $ws.Range("E5").Value = $lang.Row5.E
    $ws.Range("H5").Value = "0001-01-01 00:00:00"
    $ws.Range("I5").Value = "Include"

    # Rebuild every hyperlink on the sheet (A/B/D columns of rows 2-5) since
    # the Hyperlinks collection does not track the row inserts above.
    $ws.Hyperlinks.Delete()

    Add-ReportHyperlink $ws "A2" "https://github.com/OpenLocalizationTest/oltest/blob/f3f0bcc76a05ce10f5c256e5d9009383cc543374/e2e/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md" "0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md"
    Add-ReportHyperlink $ws "B2" "https://github.com/OpenLocalizationTest/oltest/blob/f3f0bcc76a05ce10f5c256e5d9009383cc543374/e2e/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.md" ".md"
    Add-ReportHyperlink $ws "D2" ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/086b6d37dc6c43a7cf9719a91987f1bcaf05af36/ol-handoff/OpenLocalizationTestOrg/oltest." + $suffix + "/ci/ht/0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74." + $suffix + ".xlf") ("0c05a52b-0c83-4586-9214-fbdd36cbc2b1.5ecd15e3fdae5f448c36f40c571c9b41a8f26c74." + $suffix + ".xlf")

    Add-ReportHyperlink $ws "A3" "https://github.com/OpenLocalizationTest/oltest/blob/ed81bff2ae9620a8cb3ed916d5d61fbfca77ec35/e2e/6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.md" "6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.md"
    Add-ReportHyperlink $ws "B3" "https://github.com/OpenLocalizationTest/oltest/blob/ed81bff2ae9620a8cb3ed916d5d61fbfca77ec35/e2e/6dbee821-dc55-4aa6-b9f8-a85a8a7c608a.md" ".md"
    Add-ReportHyperlink $ws "D3" ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ed81bff2ae9620a8cb3ed916d5d61fbfca77ec35/ol-handoff/OpenLocalizationTestOrg/oltest." + $suffix + "/ci/ht/" + $lang.Row3.D) $lang.Row3.D

    Add-ReportHyperlink $ws "A4" "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/8d18e275-99a5-4481-83b7-9dcffb239eac.md" "8d18e275-99a5-4481-83b7-9dcffb239eac.md"
    Add-ReportHyperlink $ws "B4" "https://github.com/OpenLocalizationTest/oltest/blob/6338203bb1567c2a4fe6ccf46ca35e8195a531b8/e2e/8d18e275-99a5-4481-83b7-9dcffb239eac.md" ".md"
    Add-ReportHyperlink $ws "D4" ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cc6d4387fd9d857c42c70350a6cb9a546af6e3d/ol-handoff/OpenLocalizationTestOrg/oltest." + $suffix + "/ci/ht/" + $lang.Row4.D) $lang.Row4.D

    Add-ReportHyperlink $ws "A5" "https://github.com/OpenLocalizationTest/oltest/blob/b87ec7a0b392ebf3a46797c41e6174f22f049a84/e2e/d1c0fc93-1653-4056-9649-c442a3e88578.md" "d1c0fc93-1653-4056-9649-c442a3e88578.md"
    Add-ReportHyperlink $ws "B5" "https://github.com/OpenLocalizationTest/oltest/blob/b87ec7a0b392ebf3a46797c41e6174f22f049a84/e2e/d1c0fc93-1653-4056-9649-c442a3e88578.md" ".md"
    Add-ReportHyperlink $ws "D5" ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b87ec7a0b392ebf3a46797c41e6174f22f049a84/ol-handoff/OpenLocalizationTestOrg/oltest." + $suffix + "/ci/ht/" + $lang.Row5.D) $lang.Row5.D
}
